# Append 7 new daily rows (2021-11-01 .. 2021-11-07) to the covid_totals
# sheet, continuing directly after the existing last row (446).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 447
$endRow = 453

$dates  = @("2021-11-01","2021-11-02","2021-11-03","2021-11-04","2021-11-05","2021-11-06","2021-11-07")
$areaType = "overview"
$areaCode = "K02000001"
$areaName = "United Kingdom"
$cumCases = @(9097311, 9130857, 9171660, 9208219, 9241916, 9272066, 9301909)
$newCases = @(40077, 33865, 41299, 37269, 34029, 30693, 30305)
$newDeaths = @(40, 293, 217, 214, 193, 155, 62)
$cumDeaths = @(140672, 140964, 141181, 141395, 141588, 141743, 141805)

# Column A holds a date-looking string (e.g. "2021-11-01"); the existing
# rows store it as literal text, so force the new cells to Text format
# first -- otherwise Excel auto-converts the string into a date serial
# number on assignment.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt ($endRow - $startRow + 1); $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $areaType
    $ws.Cells.Item($r, 3).Value = $areaCode
    $ws.Cells.Item($r, 4).Value = $areaName
    $ws.Cells.Item($r, 5).Value = $cumCases[$i]
    $ws.Cells.Item($r, 6).Value = $newCases[$i]
    $ws.Cells.Item($r, 7).Value = $newDeaths[$i]
    $ws.Cells.Item($r, 8).Value = $cumDeaths[$i]
}
